$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEJORAR")
$ws.Activate() | Out-Null
$ws = $wb.ActiveSheet

# Reproduce the prior selection state (A41:A44 selected, active cell A44)
# before the row is removed, matching the user's in-progress selection.
$ws.Range("A41:A44").Select() | Out-Null

# "GAG12301AR" (row 42) is the duplicate/unwanted product entry being removed
# from the catalog - deleting the whole row shifts every following row up by
# one (ISAALAMF16 -> A42, EA5310MT -> A43) and drops the old last row (A44).
$ws.Rows.Item(42).Delete()

# Fix the malformed product code "TF414" -> "TF.414".
$ws.Cells.Item(40, 1).Value = "TF.414"

# Final selection after the edit: single active cell A44 (now past the used
# range, same as in the target workbook).
$ws.Range("A44").Select() | Out-Null

# Scroll the viewport down so row 34 is the top visible row.
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
